# ---------------------------------------------------------------------------
# Weekly roster update:
#   - shorten column headers on the existing "09.09.2023" sheet
#     (Average Energy -> Energy, Active Guild War -> GW, Active Battles -> TW)
#   - give the sheet explicit column widths
#   - duplicate that layout into a brand-new "16.09.2023" sheet for the
#     following week, with the week's roster data (GW value differs: 8 vs 9)
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Column widths as they should appear in the saved file (stored OOXML
# "width" units). Excel's ColumnWidth property (characters) and the raw
# stored width differ by the fixed 5/6 padding offset, so compensate here.
$colWidths = @(16, 16, 11, 7, 8, 8, 4, 4, 22)

for ($c = 1; $c -le 9; $c++) {
    $ws1.Columns.Item($c).ColumnWidth = ($colWidths[$c - 1] - (5 / 6))
}

# Rename the three headers that got shortened.
$ws1.Range("F1").Value = "Energy"
$ws1.Range("G1").Value = "GW"
$ws1.Range("H1").Value = "TW"

# ---------------------------------------------------------------------------
# Add the new week's sheet right after the existing one.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "16.09.2023"

# Match the outline + page-margin presentation of the first sheet.
$ws2.Outline.SummaryRow    = 1
$ws2.Outline.SummaryColumn = 1

$ws2.PageSetup.LeftMargin   = 54
$ws2.PageSetup.RightMargin  = 54
$ws2.PageSetup.TopMargin    = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

for ($c = 1; $c -le 9; $c++) {
    $ws2.Columns.Item($c).ColumnWidth = ($colWidths[$c - 1] - (5 / 6))
}

# Header row, already using the shortened labels.
$headers = @("Player Name", "Galactic Power", "Player ID", "Level", "Role", "Energy", "GW", "TW", "Plan")
for ($c = 1; $c -le 9; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Row 2 roster data for the new week.
$ws2.Range("A2").Value = "Achernarkh Sun"
$ws2.Range("B2").Value = 5289820
$ws2.Range("C2").NumberFormat = "@"          # keep the player ID as text, not a number
$ws2.Range("C2").Value = "481289748"
$ws2.Range("D2").Value = 85
$ws2.Range("E2").Value = "Member"
$ws2.Range("F2").Value = 600
$ws2.Range("G2").Value = 8
$ws2.Range("H2").Value = 9
$ws2.Range("I2").Value = "фениксы, доктор афра"

Write-Host "Added sheet '16.09.2023' and updated headers on '09.09.2023'."
